$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.470.31"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.852.23"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "240.67"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "0.6285"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.07671"
$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("D9").Value = "0.2924"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").Value = "24.86"
$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("D11").Value = "0.07762"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").Value = "1.844.58"
$ws.Range("E12").Value = "  -2.95%  "

$ws.Range("D13").Value = "5.042"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").Value = "0.6830"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "0.00001060"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("D16").Value = "83.70"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.111.85"
$ws.Range("E17").Value = "  -2.64%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "6.200"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "29.499.61"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "229.32"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "12.37"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "7.497"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "157.78"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.1383"
$ws.Range("E26").Value = "  -1.30%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "8.429"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "17.80"
$ws.Range("E28").Value = "  +0.94%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.362"
$ws.Range("E29").Value = "  +4.98%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.465"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.05631"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.134"
$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.058"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.843"
$ws.Range("E34").Value = "  -1.04%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.165"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7034"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.601"
$ws.Range("E37").Value = "  +0.48%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.231.61"
$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01802"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.759"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.454"
$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.9039"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").Value = "102.00"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "66.32"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").Value = "7.182"
$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("D47").Value = "0.00000000118"
$ws.Range("E47").Value = "  -1.59%  "

$ws.Range("D48").Value = "0.4026"
$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D49").Value = "9.072"
$ws.Range("E49").Value = "  +1.37%  "

$ws.Range("D50").Value = "0.1159"
$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("D51").Value = "1.679"
$ws.Range("E51").Value = "  -0.91%  "
